# Regenerate save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals. This updates column G ("K") values for rows 2-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 2
    27 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
